# Update loading_percent values for the 380 kV case (rows 2-25, cols C-G,I,M,N)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 4.677597353611008
$ws.Range("D2").Value = 7.68875235114345
$ws.Range("E2").Value = 10.16783314414728
$ws.Range("F2").Value = 45.73989080449947
$ws.Range("G2").Value = 3.612113617510611
$ws.Range("I2").Value = 35.36832391276145
$ws.Range("M2").Value = 28.01105371695652
$ws.Range("N2").Value = 16.42988545052622

# Row 3
$ws.Range("C3").Value = 4.606229772817682
$ws.Range("D3").Value = 7.524750119889177
$ws.Range("E3").Value = 9.713450916160024
$ws.Range("F3").Value = 44.14010454934889
$ws.Range("G3").Value = 3.621599735080276
$ws.Range("I3").Value = 34.28672471946753
$ws.Range("M3").Value = 26.85581400791083
$ws.Range("N3").Value = 16.27876040176692

# Row 4
$ws.Range("C4").Value = 4.564646591805882
$ws.Range("D4").Value = 7.426791516403399
$ws.Range("E4").Value = 9.422364095997354
$ws.Range("F4").Value = 43.15394514143656
$ws.Range("G4").Value = 3.627686561525118
$ws.Range("I4").Value = 33.62488773166671
$ws.Range("M4").Value = 26.12835585874916
$ws.Range("N4").Value = 16.18766417315888

# Row 5
$ws.Range("C5").Value = 4.548277473316952
$ws.Range("D5").Value = 7.387612249035424
$ws.Range("E5").Value = 9.300817655920302
$ws.Range("F5").Value = 42.75178906880644
$ws.Range("G5").Value = 3.63023350432347
$ws.Range("I5").Value = 33.35619714780157
$ws.Range("M5").Value = 25.82779227178067
$ws.Range("N5").Value = 16.15100552222113

# Row 6
$ws.Range("C6").Value = 4.545594573661609
$ws.Range("D6").Value = 7.381152684156343
$ws.Range("E6").Value = 9.280461938983663
$ws.Range("F6").Value = 42.6850149334709
$ws.Range("G6").Value = 3.630660455359008
$ws.Range("I6").Value = 33.31165575458866
$ws.Range("M6").Value = 25.77764910394347
$ws.Range("N6").Value = 16.14494747966462

# Row 7
$ws.Range("C7").Value = 4.564423481482605
$ws.Range("D7").Value = 7.426260070411374
$ws.Range("E7").Value = 9.420736561274188
$ws.Range("F7").Value = 43.14852174825934
$ws.Range("G7").Value = 3.627720640457432
$ws.Range("I7").Value = 33.62125937568173
$ws.Range("M7").Value = 26.12431840657526
$ws.Range("N7").Value = 16.18716785466085

# Row 8
$ws.Range("C8").Value = 4.65253446241748
$ws.Range("D8").Value = 7.631666767706018
$ws.Range("E8").Value = 10.01373596848045
$ws.Range("F8").Value = 45.18948965783453
$ws.Range("G8").Value = 3.615330322932271
$ws.Range("I8").Value = 34.99517734597456
$ws.Range("M8").Value = 27.6167405737726
$ws.Range("N8").Value = 16.37744837869434

# Row 9
$ws.Range("C9").Value = 4.842373772579069
$ws.Range("D9").Value = 8.053916007398456
$ws.Range("E9").Value = 11.07643539354131
$ws.Range("F9").Value = 49.13320506520541
$ws.Range("G9").Value = 3.593088814300747
$ws.Range("I9").Value = 37.68904866262506
$ws.Range("M9").Value = 30.3826287433234
$ws.Range("N9").Value = 16.76236906561652

# Row 10
$ws.Range("C10").Value = 4.991225119916773
$ws.Range("D10").Value = 8.372927113883595
$ws.Range("E10").Value = 11.79153120623808
$ws.Range("F10").Value = 51.9614733099279
$ws.Range("G10").Value = 3.577965110024991
$ws.Range("I10").Value = 39.64524019996906
$ws.Range("M10").Value = 32.29756000234283
$ws.Range("N10").Value = 17.05010170269012

# Row 11
$ws.Range("C11").Value = 5.060742400457632
$ws.Range("D11").Value = 8.519327628435162
$ws.Range("E11").Value = 12.1018663754181
$ws.Range("F11").Value = 53.22714554536449
$ws.Range("G11").Value = 3.571341210393261
$ws.Range("I11").Value = 40.52594764776929
$ws.Range("M11").Value = 33.14009970385589
$ws.Range("N11").Value = 17.18159172904657

# Row 12
$ws.Range("C12").Value = 5.087305902042242
$ws.Range("D12").Value = 8.574899501854432
$ws.Range("E12").Value = 12.21718351683802
$ws.Range("F12").Value = 53.70299178736866
$ws.Range("G12").Value = 3.56886905940958
$ws.Range("I12").Value = 40.85782048778665
$ws.Range("M12").Value = 33.45482600384576
$ws.Range("N12").Value = 17.23143005434027

# Row 13
$ws.Range("C13").Value = 5.08157468284911
$ws.Range("D13").Value = 8.562925965197689
$ws.Range("E13").Value = 12.19244642744383
$ws.Range("F13").Value = 53.60066886503036
$ws.Range("G13").Value = 3.569399882602299
$ws.Range("I13").Value = 40.78642294895517
$ws.Range("M13").Value = 33.38723951403406
$ws.Range("N13").Value = 17.22069507189269

# Row 14
$ws.Range("C14").Value = 5.062923127339245
$ws.Range("D14").Value = 8.523897173468391
$ws.Range("E14").Value = 12.11139792237141
$ws.Range("F14").Value = 53.26636483317075
$ws.Range("G14").Value = 3.571137103685611
$ws.Range("I14").Value = 40.55328540850653
$ws.Range("M14").Value = 33.16608022820268
$ws.Range("N14").Value = 17.18569126642279

# Row 15
$ws.Range("C15").Value = 5.051529000889976
$ws.Range("D15").Value = 8.50000678653052
$ws.Range("E15").Value = 12.06146558545316
$ws.Range("F15").Value = 53.06113465757269
$ws.Range("G15").Value = 3.572205894006631
$ws.Range("I15").Value = 40.41026067827924
$ws.Range("M15").Value = 33.03004467033102
$ws.Range("N15").Value = 17.16425521899438

# Row 16
$ws.Range("C16").Value = 4.986716543428313
$ws.Range("D16").Value = 8.363380850336904
$ws.Range("E16").Value = 11.77094461937706
$ws.Range("F16").Value = 51.87829985765489
$ws.Range("G16").Value = 3.57840309436571
$ws.Range("I16").Value = 39.5874718150685
$ws.Range("M16").Value = 32.24190380025227
$ws.Range("N16").Value = 17.04151712216982

# Row 17
$ws.Range("C17").Value = 4.94740361910584
$ws.Range("D17").Value = 8.279856218291583
$ws.Range("E17").Value = 11.58884934470143
$ws.Range("F17").Value = 51.14698900358683
$ws.Range("G17").Value = 3.582269987964083
$ws.Range("I17").Value = 39.08013497182735
$ws.Range("M17").Value = 31.75092164530349
$ws.Range("N17").Value = 16.96634590243813

# Row 18
$ws.Range("C18").Value = 4.924962917177181
$ws.Range("D18").Value = 8.231938517388771
$ws.Range("E18").Value = 11.48270659338423
$ws.Range("F18").Value = 50.72441124133292
$ws.Range("G18").Value = 3.584518252582381
$ws.Range("I18").Value = 38.7874813841529
$ws.Range("M18").Value = 31.46584295489365
$ws.Range("N18").Value = 16.92316816410098

# Row 19
$ws.Range("C19").Value = 4.917394886692421
$ws.Range("D19").Value = 8.215737206655714
$ws.Range("E19").Value = 11.44652854764335
$ws.Range("F19").Value = 50.58101317691084
$ws.Range("G19").Value = 3.585283639569809
$ws.Range("I19").Value = 38.68825852365664
$ws.Range("M19").Value = 31.36886730408951
$ws.Range("N19").Value = 16.90856025038814

# Row 20
$ws.Range("C20").Value = 4.951570992993675
$ws.Range("D20").Value = 8.288735158134584
$ws.Range("E20").Value = 11.60837962706958
$ws.Range("F20").Value = 51.22504308035647
$ws.Range("G20").Value = 3.581855857716033
$ws.Range("I20").Value = 39.13423186022216
$ws.Range("M20").Value = 31.80346645200665
$ws.Range("N20").Value = 16.97434221896517

# Row 21
$ws.Range("C21").Value = 5.068395223967595
$ws.Range("D21").Value = 8.535357654649754
$ws.Range("E21").Value = 12.13526388131278
$ws.Range("F21").Value = 53.36465443855377
$ws.Range("G21").Value = 3.570625863022913
$ws.Range("I21").Value = 40.62181017734481
$ws.Range("M21").Value = 33.2311590123585
$ws.Range("N21").Value = 17.19597180442413

# Row 22
$ws.Range("C22").Value = 5.146129149037484
$ws.Range("D22").Value = 8.697294275360649
$ws.Range("E22").Value = 12.46677325482046
$ws.Range("F22").Value = 54.74282309321102
$ws.Range("G22").Value = 3.563496970760956
$ws.Range("I22").Value = 41.58439417687519
$ws.Range("M22").Value = 34.13895103719984
$ws.Range("N22").Value = 17.34106891038293

# Row 23
$ws.Range("C23").Value = 5.104521371439909
$ws.Range("D23").Value = 8.610812752639205
$ws.Range("E23").Value = 12.2910289828882
$ws.Range("F23").Value = 54.00924370812101
$ws.Range("G23").Value = 3.567282740572345
$ws.Range("I23").Value = 41.07162145076565
$ws.Range("M23").Value = 33.65682202394132
$ws.Range("N23").Value = 17.26361817016779

# Row 24
$ws.Range("C24").Value = 4.949686420770256
$ws.Range("D24").Value = 8.284720668687125
$ws.Range("E24").Value = 11.59955451078944
$ws.Range("F24").Value = 51.18976145674738
$ws.Range("G24").Value = 3.582043007661242
$ws.Range("I24").Value = 39.10977768187004
$ws.Range("M24").Value = 31.77971966975065
$ws.Range("N24").Value = 16.97072695911588

# Row 25
$ws.Range("C25").Value = 4.789292614536521
$ws.Range("D25").Value = 7.937923406987789
$ws.Range("E25").Value = 10.80019746726393
$ws.Range("F25").Value = 48.07626571390547
$ws.Range("G25").Value = 3.598889332182854
$ws.Range("I25").Value = 36.9627857486348
$ws.Range("M25").Value = 29.65374666888442
$ws.Range("N25").Value = 16.65720871355514

